# Generate Report for Archive
#
# This localization-status workbook is produced by a CI report generator;
# this run refreshes two things that changed between generations:
#
#   1. The two localization targets (zh-cn / de-de) for the first source
#      file have moved on in the pipeline: their status flips from
#      "Ready for handoff" to "In Translation" everywhere it is shown
#      (Overview!E2:F3, and the Status column -- column C -- on the
#      per-locale "zh-cn" / "de-de" detail sheets).
#   2. The "Status" column is narrower in the refreshed report, so its
#      width shrinks from ~17.22 to ~13.41 "raw" Excel width units on all
#      three sheets.
#
# Helper: Excel's ColumnWidth (character units) is silently re-quantized
# to whole pixels of the Normal-style font before it is stored as the
# worksheet's raw <col width="..."> (character) units, so we can't just
# assign the target raw width. Back it out: pick the integer pixel count
# that lands closest to the desired raw width, then convert that back to
# the character-based ColumnWidth Excel expects as input.
function Set-RawColumnWidth($Column, $TargetRawWidth) {
    $mdw = 6.0
    $exactPixels = ($TargetRawWidth * $mdw) - 5.0
    $lo = [Math]::Floor($exactPixels)
    $hi = [Math]::Ceiling($exactPixels)
    $bestErr = 1000000
    $bestChars = 0
    foreach ($p in @($lo, $hi)) {
        $stored = ($p + 5.0) / $mdw
        $err = [Math]::Abs($stored - $TargetRawWidth)
        if ($err -lt $bestErr) {
            $bestErr = $err
            $bestChars = $p / $mdw
        }
    }
    $Column.ColumnWidth = $bestChars
}

$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# --- 2. Status column width: ~17.2159881591797 -> ~13.4101845877511 ---
Set-RawColumnWidth $overview.Columns.Item(5) 13.4101845877511
Set-RawColumnWidth $overview.Columns.Item(6) 13.4101845877511
Set-RawColumnWidth $zhcn.Columns.Item(3) 13.4101845877511
Set-RawColumnWidth $dede.Columns.Item(3) 13.4101845877511
